$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2022-10-29"

# Update the October label (row 11, column A) to the new "through" date
$ws.Range("A11").Value = "October (through 10-29)"

# Update January row (row 2), 2022 column (I) value
$ws.Range("I2").Value = 162

# Update October row (row 11) values for all years
$ws.Range("B11").Value = 29
$ws.Range("C11").Value = 53
$ws.Range("D11").Value = 78
$ws.Range("E11").Value = 60
$ws.Range("F11").Value = 57
$ws.Range("G11").Value = 139
$ws.Range("H11").Value = 178
$ws.Range("I11").Value = 114

# Update Total row (row 12) values for all years
$ws.Range("B12").Value = 255
$ws.Range("C12").Value = 482
$ws.Range("D12").Value = 705
$ws.Range("E12").Value = 608
$ws.Range("F12").Value = 479
$ws.Range("G12").Value = 1040
$ws.Range("H12").Value = 1425
$ws.Range("I12").Value = 1391
